$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "47+48=95"
$t.Cell(1,2).Range.Text = "90-67=23"
$t.Cell(1,3).Range.Text = "32+37=69"
$t.Cell(1,4).Range.Text = "20+9=29"
$t.Cell(1,5).Range.Text = "74-50=24"
$t.Cell(2,1).Range.Text = "86-5=81"
$t.Cell(2,2).Range.Text = "34+49=83"
$t.Cell(2,3).Range.Text = "38+39=77"
$t.Cell(2,4).Range.Text = "29+19=48"
$t.Cell(2,5).Range.Text = "82-77=5"
$t.Cell(3,1).Range.Text = "26+59=85"
$t.Cell(3,2).Range.Text = "2+96=98"
$t.Cell(3,3).Range.Text = "80-3=77"
$t.Cell(3,4).Range.Text = "51+3=54"
$t.Cell(3,5).Range.Text = "36+33=69"
$t.Cell(4,1).Range.Text = "87-41=46"
$t.Cell(4,2).Range.Text = "6+32=38"
$t.Cell(4,3).Range.Text = "29-8=21"
$t.Cell(4,4).Range.Text = "66-30=36"
$t.Cell(4,5).Range.Text = "28+49=77"
$t.Cell(5,1).Range.Text = "77-53=24"
$t.Cell(5,2).Range.Text = "90-32=58"
$t.Cell(5,3).Range.Text = "70+26=96"
$t.Cell(5,4).Range.Text = "5+2=7"
$t.Cell(5,5).Range.Text = "23+6=29"
$t.Cell(6,1).Range.Text = "43-11=32"
$t.Cell(6,2).Range.Text = "66+30=96"
$t.Cell(6,3).Range.Text = "88-55=33"
$t.Cell(6,4).Range.Text = "31-22=9"
$t.Cell(6,5).Range.Text = "99-79=20"
$t.Cell(7,1).Range.Text = "8+56=64"
$t.Cell(7,2).Range.Text = "75-65=10"
$t.Cell(7,3).Range.Text = "85+1=86"
$t.Cell(7,4).Range.Text = "81-76=5"
$t.Cell(7,5).Range.Text = "25+16=41"
$t.Cell(8,1).Range.Text = "32+47=79"
$t.Cell(8,2).Range.Text = "83+10=93"
$t.Cell(8,3).Range.Text = "76-24=52"
$t.Cell(8,4).Range.Text = "51-47=4"
$t.Cell(8,5).Range.Text = "56+41=97"
$t.Cell(9,1).Range.Text = "67-13=54"
$t.Cell(9,2).Range.Text = "1+82=83"
$t.Cell(9,3).Range.Text = "73-19=54"
$t.Cell(9,4).Range.Text = "9+24=33"
$t.Cell(9,5).Range.Text = "7+12=19"
$t.Cell(10,1).Range.Text = "8+63=71"
$t.Cell(10,2).Range.Text = "65-59=6"
$t.Cell(10,3).Range.Text = "57+35=92"
$t.Cell(10,4).Range.Text = "24+12=36"
$t.Cell(10,5).Range.Text = "60-18=42"
$t.Cell(11,1).Range.Text = "9+14=23"
$t.Cell(11,2).Range.Text = "28+43=71"
$t.Cell(11,3).Range.Text = "17-4=13"
$t.Cell(11,4).Range.Text = "21+70=91"
$t.Cell(11,5).Range.Text = "17-1=16"
$t.Cell(12,1).Range.Text = "82-44=38"
$t.Cell(12,2).Range.Text = "39+13=52"
$t.Cell(12,3).Range.Text = "23+26=49"
$t.Cell(12,4).Range.Text = "10+42=52"
$t.Cell(12,5).Range.Text = "64-41=23"
$t.Cell(13,1).Range.Text = "81+7=88"
$t.Cell(13,2).Range.Text = "13+5=18"
$t.Cell(13,3).Range.Text = "59-1=58"
$t.Cell(13,4).Range.Text = "65-57=8"
$t.Cell(13,5).Range.Text = "27-21=6"
$t.Cell(14,1).Range.Text = "78-30=48"
$t.Cell(14,2).Range.Text = "92-44=48"
$t.Cell(14,3).Range.Text = "70-35=35"
$t.Cell(14,4).Range.Text = "66-23=43"
$t.Cell(14,5).Range.Text = "7+57=64"
$t.Cell(15,1).Range.Text = "88-83=5"
$t.Cell(15,2).Range.Text = "91-38=53"
$t.Cell(15,3).Range.Text = "58+39=97"
$t.Cell(15,4).Range.Text = "65-21=44"
$t.Cell(15,5).Range.Text = "65+6=71"
$t.Cell(16,1).Range.Text = "80-39=41"
$t.Cell(16,2).Range.Text = "16+20=36"
$t.Cell(16,3).Range.Text = "3+86=89"
$t.Cell(16,4).Range.Text = "23+49=72"
$t.Cell(16,5).Range.Text = "15+2=17"
$t.Cell(17,1).Range.Text = "13+25=38"
$t.Cell(17,2).Range.Text = "86-64=22"
$t.Cell(17,3).Range.Text = "30+52=82"
$t.Cell(17,4).Range.Text = "88-71=17"
$t.Cell(17,5).Range.Text = "31-28=3"
$t.Cell(18,1).Range.Text = "43-43=0"
$t.Cell(18,2).Range.Text = "45+27=72"
$t.Cell(18,3).Range.Text = "90-57=33"
$t.Cell(18,4).Range.Text = "56-12=44"
$t.Cell(18,5).Range.Text = "89-51=38"
$t.Cell(19,1).Range.Text = "36+48=84"
$t.Cell(19,2).Range.Text = "94-63=31"
$t.Cell(19,3).Range.Text = "71-40=31"
$t.Cell(19,4).Range.Text = "6+77=83"
$t.Cell(19,5).Range.Text = "9+54=63"
$t.Cell(20,1).Range.Text = "28-21=7"
$t.Cell(20,2).Range.Text = "75-53=22"
$t.Cell(20,3).Range.Text = "64+9=73"
$t.Cell(20,4).Range.Text = "5+70=75"
$t.Cell(20,5).Range.Text = "41-33=8"
